# Add a new "Thank you" closing slide (Title and Content layout) as the
# last slide of the deck.

$p = $ppt.ActivePresentation

# "Заголовок и объект" (Title and Content) is CustomLayout index 2 on the
# (single) slide master used throughout this deck.
$layout = $p.SlideMaster.CustomLayouts.Item(2)

$newIndex = $p.Slides.Count + 1
$s = $p.Slides.AddSlide($newIndex, $layout)

# Match the locale-default placeholder names used by the rest of the deck.
$s.Shapes.Item(1).Name = "Заголовок 1"
$s.Shapes.Item(2).Name = "Объект 2"

# --- Title placeholder ------------------------------------------------
$title = $s.Shapes.Item(1).TextFrame.TextRange
$title.Text = "`t`t`t`tThank "
$titleRun2 = $title.InsertAfter("you for attention")

# --- Body / content placeholder ---------------------------------------
$body = $s.Shapes.Item(2).TextFrame.TextRange
$body.Text = "Performed "
$run2 = $body.InsertAfter("by ")
$run3 = $run2.InsertAfter("Oktamov")
$run4 = $run3.InsertAfter(" ")
$run5 = $run4.InsertAfter("Saidakmal")

$full = $s.Shapes.Item(2).TextFrame.TextRange
$full.Font.Size = 40
$full.Font.Bold = $true
$full.Font.Italic = $true

# No bullet, no left margin/indent for this single paragraph.
$full.ParagraphFormat.Bullet.Visible = $false
$level1 = $s.Shapes.Item(2).TextFrame.Ruler.Levels.Item(1)
$level1.FirstMargin = 0
$level1.LeftMargin = 0
